# Auto update Excel log
# Appends newly-logged sensor readings (2026-01-28, ~14:57-14:58) to the
# PIR, Humidity and Temperature sheets. NumberFormat is forced to "@"
# (Text) on columns whose literal values look like dates/percentages so
# Excel doesn't silently re-type them as numbers/dates on input.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")
$PIRData = @(
    @("2026-01-28","14:57:45","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:57:49","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:57:54","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:57:59","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:58:04","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:58:09","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:58:14","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:58:19","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:58:24","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:58:29","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:58:34","14:00","Bathroom","No Motion","Inactive"),
    @("2026-01-28","14:58:39","14:00","Bathroom","No Motion","Inactive")
)
$PIRTextCols = @(1)
$startRow = 123
for ($i = 0; $i -lt $PIRData.Length; $i++) {
    $rowData = $PIRData[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $cell = $ws.Cells.Item($r, $j + 1)
        if ($PIRTextCols -contains ($j + 1)) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $rowData[$j]
    }
}

$ws = $wb.Worksheets.Item("Humidity")
$HumidityData = @(
    @("2026-01-28","14:57:44","14:00","Bathroom","87.0%","Active"),
    @("2026-01-28","14:57:45","14:00","Bathroom","88.0%","Active"),
    @("2026-01-28","14:57:49","14:00","Bathroom","87.9%","Active"),
    @("2026-01-28","14:57:53","14:00","Bathroom","87.0%","Active"),
    @("2026-01-28","14:57:57","14:00","Bathroom","87.9%","Active"),
    @("2026-01-28","14:58:01","14:00","Bathroom","87.0%","Active"),
    @("2026-01-28","14:58:05","14:00","Bathroom","87.9%","Active"),
    @("2026-01-28","14:58:09","14:00","Bathroom","87.9%","Active"),
    @("2026-01-28","14:58:13","14:00","Bathroom","87.1%","Active"),
    @("2026-01-28","14:58:17","14:00","Bathroom","88.0%","Active"),
    @("2026-01-28","14:58:25","14:00","Bathroom","87.1%","Active"),
    @("2026-01-28","14:58:29","14:00","Bathroom","88.0%","Active"),
    @("2026-01-28","14:58:33","14:00","Bathroom","87.1%","Active"),
    @("2026-01-28","14:58:38","14:00","Bathroom","88.1%","Active")
)
$HumidityTextCols = @(1,5)
$startRow = 117
for ($i = 0; $i -lt $HumidityData.Length; $i++) {
    $rowData = $HumidityData[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $cell = $ws.Cells.Item($r, $j + 1)
        if ($HumidityTextCols -contains ($j + 1)) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $rowData[$j]
    }
}

$ws = $wb.Worksheets.Item("Temperature")
$TemperatureData = @(
    @("2026-01-28","14:57:44","14:00","Bathroom","22.9C","Active"),
    @("2026-01-28","14:57:46","14:00","Bathroom","22.9C","Active"),
    @("2026-01-28","14:57:50","14:00","Bathroom","22.9C","Active"),
    @("2026-01-28","14:57:54","14:00","Bathroom","22.9C","Active"),
    @("2026-01-28","14:57:58","14:00","Bathroom","22.9C","Active"),
    @("2026-01-28","14:58:02","14:00","Bathroom","22.9C","Active"),
    @("2026-01-28","14:58:06","14:00","Bathroom","22.9C","Active"),
    @("2026-01-28","14:58:10","14:00","Bathroom","22.9C","Active"),
    @("2026-01-28","14:58:14","14:00","Bathroom","22.9C","Active"),
    @("2026-01-28","14:58:18","14:00","Bathroom","22.9C","Active"),
    @("2026-01-28","14:58:26","14:00","Bathroom","22.9C","Active"),
    @("2026-01-28","14:58:30","14:00","Bathroom","22.9C","Active"),
    @("2026-01-28","14:58:34","14:00","Bathroom","22.9C","Active"),
    @("2026-01-28","14:58:38","14:00","Bathroom","22.9C","Active")
)
$TemperatureTextCols = @(1)
$startRow = 117
for ($i = 0; $i -lt $TemperatureData.Length; $i++) {
    $rowData = $TemperatureData[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $cell = $ws.Cells.Item($r, $j + 1)
        if ($TemperatureTextCols -contains ($j + 1)) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $rowData[$j]
    }
}
